$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15, pushing the existing rows 15-21 down to 16-22.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with this week's data (carries D15's date style
# forward automatically from the row-insert, matching the other date cells).
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 44511
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 100112026
$ws.Range("G15").Value = "Haba"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 7500
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Provincia de Diguillín"
$ws.Range("P15").Value = 300
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
